# Weekly cryptos.xlsx refresh (GitHub Actions bot) -- updates the Price (D) and
# Volume(1h) (E) columns for every coin row (2-51) with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price strings look like plain decimals (e.g. "1.700",
# "1.003", "0.2658") that Excel would otherwise silently coerce into numbers,
# which would strip meaningful trailing zeros. Force just those specific cells
# to Text format before writing so the literal digits are preserved verbatim,
# exactly like the already-textual prices elsewhere in the column.
$textPriceRows = @(4,5,7,8,9,10,11,12,15,17,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,41,43,44,47,48,49,50,51)
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '27.275.96'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Value = '1.703.42'
$ws.Range("E3").Value = '  -1.16%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '223.40'
$ws.Range("E5").Value = '  -1.07%  '
$ws.Range("E6").Value = '  -1.25%  '
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '0.2658'
$ws.Range("E8").Value = '  -0.89%  '
$ws.Range("D9").Value = '0.06579'
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("D10").Value = '20.71'
$ws.Range("E10").Value = '  -4.25%  '
$ws.Range("D11").Value = '0.07621'
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").Value = '4.490'
$ws.Range("E12").Value = '  -3.12%  '
$ws.Range("D13").Value = '1.707.46'
$ws.Range("E13").Value = '  -1.06%  '
$ws.Range("D14").Value = '1.940.30'
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("D15").Value = '0.5777'
$ws.Range("E15").Value = '  -1.68%  '
$ws.Range("D16").Value = '0.0₅8141'
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("D17").Value = '67.47'
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = '27.289.28'
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("D19").Value = '215.24'
$ws.Range("E19").Value = '  -3.71%  '
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = '4.606'
$ws.Range("E21").Value = '  -2.74%  '
$ws.Range("D22").Value = '10.35'
$ws.Range("E22").Value = '  -3.28%  '
$ws.Range("D23").Value = '5.959'
$ws.Range("E23").Value = '  -2.33%  '
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").Value = '143.89'
$ws.Range("E25").Value = '  -2.89%  '
$ws.Range("D26").Value = '1.700'
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("D27").Value = '0.1198'
$ws.Range("E27").Value = '  -2.88%  '
$ws.Range("D28").Value = '7.199'
$ws.Range("E28").Value = '  -2.85%  '
$ws.Range("D29").Value = '16.11'
$ws.Range("E29").Value = '  -3.43%  '
$ws.Range("D30").Value = '0.05361'
$ws.Range("E30").Value = '  -3.72%  '
$ws.Range("D31").Value = '1.285'
$ws.Range("E31").Value = '  -1.44%  '
$ws.Range("D32").Value = '3.460'
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("D33").Value = '3.399'
$ws.Range("E33").Value = '  -1.95%  '
$ws.Range("D34").Value = '1.639'
$ws.Range("E34").Value = '  -1.46%  '
$ws.Range("D35").Value = '2.863'
$ws.Range("E35").Value = '  +1.62%  '
$ws.Range("D36").Value = '2.412'
$ws.Range("E36").Value = '  -1.58%  '
$ws.Range("D37").Value = '0.9451'
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("D38").Value = '0.5800'
$ws.Range("E38").Value = '  -2.21%  '
$ws.Range("D39").Value = '0.01624'
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("D41").Value = '1.003'
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").Value = '1.039.79'
$ws.Range("E42").Value = '  -1.74%  '
$ws.Range("D43").Value = '0.8397'
$ws.Range("E43").Value = '  -1.99%  '
$ws.Range("D44").Value = '100.97'
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").Value = '1.847.52'
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("D46").Value = '0.0₈113'
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("D47").Value = '57.72'
$ws.Range("E47").Value = '  -2.17%  '
$ws.Range("D48").Value = '0.4513'
$ws.Range("E48").Value = '  +1.70%  '
$ws.Range("D49").Value = '1.005'
$ws.Range("E49").Value = '  +0.30%  '
$ws.Range("D50").Value = '8.064'
$ws.Range("E50").Value = '  -1.89%  '
$ws.Range("D51").Value = '0.05226'
$ws.Range("E51").Value = '  -1.01%  '
